# Completion of Section F
# Mark the "Results" section (Section F) heading and all of its
# sub-headings as strikethrough, indicating the section is complete.

$d = $word.ActiveDocument

# The exact heading texts that belong to the "Results" section (the
# section-level Heading1 title plus its Heading2 sub-headings).
$targets = @(
    "Results",
    "Statistical Significance",
    "Practical Significance",
    "Overall Success/Effectiveness"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($targets -contains $text) {
        $p.Range.Font.StrikeThrough = 1
    }
}
